# "excel sheet bug fixed"
# The sheet was blank (A1:A1) with no data. Populate it with the student
# admission table: a bold header row (A1:O1) plus two data rows, set
# sensible column widths, and make sure numeric-looking reference columns
# (contact number / subject marks) are stored as text while the computed
# Percentage Marks column stays numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column widths ------------------------------------------------------
$widths = @(30, 30, 30, 60, 20, 50, 50, 50, 30, 30, 30, 30, 30, 30, 30)
for ($i = 0; $i -lt $widths.Length; $i++) {
    # ColumnWidth is specified in "characters"; Excel stores the serialized
    # width with a fixed ~0.8333 char padding added on top, so back it out
    # here to land on the exact target width once saved.
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i] - (5 / 6)
}

# --- header row (row 1), bold ------------------------------------------
$headers = @(
    "Student Name",
    "Father's Name",
    "Mother Name",
    "Email",
    "Gender",
    "Contact Number",
    "Present Address",
    "Permanent Address",
    "X Marks",
    "XII Marks",
    "Board",
    "Physics Marks",
    "Chemistry Marks",
    "Maths Marks",
    "Percentage Marks"
)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $headers.Length))
$headerRange.Font.Bold = $true

# --- helper: write a value, forcing text storage for numeric-looking ----
# strings so e.g. "93" / "9073042220" don't get silently coerced to
# numbers (which would drop meaning like leading zeros / phone formatting).
# The text number-format is only needed transiently to pin the storage
# type; restoring the "Normal" style afterwards drops the leftover
# format/style index so the cell ends up plain (unstyled) text, matching
# a normal data cell.
function Set-TextCell($sheet, $row, $col, $val) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- data rows ------------------------------------------------------------
$row2 = @(
    "adityaprava sen",
    "arup sen",
    "mousumi sen",
    "adityapravasen.0911@gmail.com",
    "Male",
    "9073042220",
    "rabindra pally, bramhapur",
    "rabindra pally, bramhapur",
    "93",
    "89",
    "CBSE",
    "95",
    "75",
    "95"
)
$row3 = @(
    "jhonny doe",
    "john doe",
    "jane doe",
    "jhonny123@gmail.com",
    "Male",
    "911",
    "chicago",
    "chicago",
    "80",
    "90",
    "ISC",
    "70",
    "80",
    "90"
)
# columns (1-based) among the 14 text fields above whose values look like
# numbers and must be pinned to text storage
$forceTextCols = @(6, 9, 10, 12, 13, 14)

for ($c = 0; $c -lt $row2.Length; $c++) {
    $col = $c + 1
    if ($forceTextCols -contains $col) {
        Set-TextCell $ws 2 $col $row2[$c]
    } else {
        $ws.Cells.Item(2, $col).Value = $row2[$c]
    }
}
$ws.Cells.Item(2, 15).Value = 88

for ($c = 0; $c -lt $row3.Length; $c++) {
    $col = $c + 1
    if ($forceTextCols -contains $col) {
        Set-TextCell $ws 3 $col $row3[$c]
    } else {
        $ws.Cells.Item(3, $col).Value = $row3[$c]
    }
}
$ws.Cells.Item(3, 15).Value = 80
